$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on modified Price cells so formatted numbers (thousand-dot style, trailing zeros, etc.) are preserved exactly
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.660.41"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "3.863.18"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "602.39"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").Value = "168.69"
$ws.Range("E6").Value = "  +2.56%  "

$ws.Range("D7").Value = "3.865.10"
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("D11").Value = "6.35"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "37.67"
$ws.Range("E14").Value = "  +1.28%  "

$ws.Range("D15").Value = "4.513.61"
$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("D16").Value = "3.868.08"
$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("D17").Value = "68.786.46"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "7.57"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "18.40"
$ws.Range("E19").Value = "  +7.03%  "

$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("D21").Value = "10.81"
$ws.Range("E21").Value = "  -3.95%  "

$ws.Range("D22").Value = "478.36"
$ws.Range("E22").Value = "  -2.24%  "

$ws.Range("D23").Value = "0.739"
$ws.Range("E23").Value = "  +1.87%  "

$ws.Range("E24").Value = "  -2.77%  "

$ws.Range("E25").Value = "  +0.26%  "

$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").Value = "12.35"
$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("D31").Value = "4.015.72"
$ws.Range("E31").Value = "  -1.61%  "

$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("E33").Value = "  -3.13%  "

$ws.Range("D34").Value = "31.22"
$ws.Range("E34").Value = "  -4.18%  "

$ws.Range("D35").Value = "3.833.64"
$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("E37").Value = "  +1.19%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "1.02"
$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").Value = "3.31"
$ws.Range("E40").Value = "  +8.62%  "

$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").Value = "2.02"
$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("D44").Value = "429.10"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").Value = "47.85"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D47").Value = "8.63"
$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "0.000273"
$ws.Range("E48").Value = "  +13.14%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "142.58"
$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("D50").Value = "0.0361"
$ws.Range("E50").Value = "  +0.80%  "

$ws.Range("D51").Value = "39.26"
$ws.Range("E51").Value = "  +0.52%  "

